# Auto-generated: update leve market-data derived columns (H-N) across sheets
# to reflect refreshed Universalis price data, per scheduled runner commit.
$wb = $excel.ActiveWorkbook

# ALC row 19: Unbreak My Heart | Roof Tile
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 889012.9399999999
$ws.Range("I19").Value = 1212206.9
$ws.Range("J19").Value = 229.75
$ws.Range("K19").Value = 1212206.9
$ws.Range("L19").Value = 229.75
$ws.Range("M19").Value = -1212031.9
$ws.Range("N19").Value = -579.75

# ALC row 80: Cleansing the Wicked Humours | Hallowed Water
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H80").Value = 1225
$ws.Range("I80").Value = 200
$ws.Range("J80").Value = 1566.6666
$ws.Range("K80").Value = 600
$ws.Range("L80").Value = 4699.9998
$ws.Range("M80").Value = 398
$ws.Range("N80").Value = -6695.9998

# ALC row 83: Washing Away the Sins (L) | Hallowed Water
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H83").Value = 1225
$ws.Range("I83").Value = 200
$ws.Range("J83").Value = 1566.6666
$ws.Range("K83").Value = 1800
$ws.Range("L83").Value = 14099.9994
$ws.Range("M83").Value = 3192
$ws.Range("N83").Value = -24083.9994

# ALC row 112: Making Ends Meet | Superior Spiritbond Potion
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 10102579
$ws.Range("J112").Value = 1600.7916
$ws.Range("L112").Value = 4802.3748
$ws.Range("N112").Value = -7018.3748

# ALC row 123: Nearly Bare | Gaja Grimoire
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H123").Value = 43280
$ws.Range("J123").Value = 43280
$ws.Range("L123").Value = 43280
$ws.Range("N123").Value = -53080

# ALC row 129: Practical Command | Commanding Craftsman's Draught
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 25298.975
$ws.Range("J129").Value = 26645.432
$ws.Range("L129").Value = 79936.296
$ws.Range("N129").Value = -89936.296

# ALC row 132: Fast-forwarding Flora | Growth Formula Lambda
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 133076.1
$ws.Range("I132").Value = 193160.42
$ws.Range("J132").Value = 6899
$ws.Range("K132").Value = 579481.26
$ws.Range("L132").Value = 20697
$ws.Range("M132").Value = -576951.26
$ws.Range("N132").Value = -25757

# ALC row 137: Cutting Edge of Culinary Quality | Magnesia Whetstone
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 4143.077
$ws.Range("I137").Value = 3511.276
$ws.Range("K137").Value = 10533.828
$ws.Range("M137").Value = -7983.828

# ARM row 61: Dealing with the Tough Stuff | Cobalt Ingot
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2879.5334
$ws.Range("J61").Value = 3174.4167
$ws.Range("L61").Value = 3174.4167
$ws.Range("N61").Value = -3598.4167

# ARM row 132: Don't Bore Me, Ore Me | Mountain Chromite Ingot
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 3623.6072
$ws.Range("I132").Value = 2660.1667
$ws.Range("J132").Value = 5357.8
$ws.Range("K132").Value = 7980.500100000001
$ws.Range("L132").Value = 16073.4
$ws.Range("M132").Value = -5450.500100000001
$ws.Range("N132").Value = -21133.4

# ARM row 136: Metal with Mettle | Cobalt Tungsten Ingot
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 2879.5334
$ws.Range("J136").Value = 3174.4167
$ws.Range("L136").Value = 9523.250100000001
$ws.Range("N136").Value = -14623.2501

# ARM row 137: Odd Instruments | Cobalt Tungsten Alembic
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H137").Value = 51549
$ws.Range("J137").Value = 51549
$ws.Range("L137").Value = 51549
$ws.Range("N137").Value = -61749

# BSM row 12: A Hit Job | Bronze Chaser Hammer
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H12").Value = 2250
$ws.Range("I12").Value = 500
$ws.Range("J12").Value = 7500
$ws.Range("K12").Value = 500
$ws.Range("L12").Value = 7500
$ws.Range("M12").Value = -332
$ws.Range("N12").Value = -7836

# BSM row 107: The Gold Experience | Deepgold Nugget
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 1032.8572
$ws.Range("I107").Value = 1032.8572
$ws.Range("K107").Value = 1032.8572
$ws.Range("M107").Value = 887.1428000000001

# BSM row 134: Ruthenium Supremium | Ruthenium Ingot
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 4041.697
$ws.Range("I134").Value = 2205.8096
$ws.Range("J134").Value = 7254.5
$ws.Range("K134").Value = 6617.4288
$ws.Range("L134").Value = 21763.5
$ws.Range("M134").Value = -4082.4288
$ws.Range("N134").Value = -26833.5

# CRP row 4: A Clogful of Camaraderie | Maple Clogs
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 4000
$ws.Range("I4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("M4").Value = $null

# CRP row 31: Wall Not Found | Walnut Lumber
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2830.449
$ws.Range("I31").Value = 1196.5
$ws.Range("K31").Value = 1196.5
$ws.Range("M31").Value = -901.5

# CRP row 34: Armoires of the Rich and Famous | Walnut Lumber
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 2830.449
$ws.Range("I34").Value = 1196.5
$ws.Range("K34").Value = 1196.5
$ws.Range("M34").Value = -994.5

# CRP row 132: Hull Lotta Damage | Ginseng Lumber
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 4560.5
$ws.Range("I132").Value = 4233.5386
$ws.Range("J132").Value = 4843.8667
$ws.Range("K132").Value = 12700.6158
$ws.Range("L132").Value = 14531.6001
$ws.Range("M132").Value = -10170.6158
$ws.Range("N132").Value = -19591.6001

# CRP row 134: Wood You Be Quiet | Ceiba Lumber
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 6447.75
$ws.Range("I134").Value = 7371.5
$ws.Range("J134").Value = 4600.25
$ws.Range("K134").Value = 22114.5
$ws.Range("L134").Value = 13800.75
$ws.Range("M134").Value = -19579.5
$ws.Range("N134").Value = -18870.75

# CRP row 138: Bow Out | Acacia Longbow
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H138").Value = 35012.637
$ws.Range("J138").Value = 35012.637
$ws.Range("L138").Value = 35012.637
$ws.Range("N138").Value = -45292.637

# CUL row 113: Can't Eat Just One | Night Vinegar
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 528.51514
$ws.Range("I113").Value = 532.0476
$ws.Range("K113").Value = 1596.1428
$ws.Range("M113").Value = 573.8571999999999

# CUL row 124: Bobbing for Compliments | Island Miq'abob
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H124").Value = 3150
$ws.Range("I124").Value = 966.6667
$ws.Range("J124").Value = 5333.3335
$ws.Range("K124").Value = 2900.0001
$ws.Range("L124").Value = 16000.0005
$ws.Range("M124").Value = 2009.9999
$ws.Range("N124").Value = -25820.0005

# CUL row 131: The Mountain Steeped | Tsai tou Vounou
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 10000972
$ws.Range("J131").Value = 828.7778
$ws.Range("L131").Value = 2486.3334
$ws.Range("N131").Value = -12566.3334

# GSM row 4: Arms for the Poor | Bone Brand
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H4").Value = 28999.166
$ws.Range("J4").Value = 28999.166
$ws.Range("L4").Value = 28999.166
$ws.Range("N4").Value = -29223.166

# GSM row 70: Sky Is the Limit | Mythrite Ingot
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6481.788
$ws.Range("I70").Value = 5773.0386
$ws.Range("K70").Value = 5773.0386
$ws.Range("M70").Value = -5503.0386

# GSM row 73: Hulls of Broken Dreams (L) | Mythrite Ingot
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 6481.788
$ws.Range("I73").Value = 5773.0386
$ws.Range("K73").Value = 5773.0386
$ws.Range("M73").Value = -4837.0386

# GSM row 102: Put the Metal to the Peddle | Durium Ingot
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 3352.0908
$ws.Range("I102").Value = 2629.0833
$ws.Range("J102").Value = 4219.7
$ws.Range("K102").Value = 2629.0833
$ws.Range("L102").Value = 4219.7
$ws.Range("M102").Value = -1007.0833
$ws.Range("N102").Value = -7463.7

# GSM row 113: Copious Crystal Cannons | Manasilver Nugget
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 1203.2
$ws.Range("I113").Value = 1182
$ws.Range("J113").Value = 1500
$ws.Range("K113").Value = 1182
$ws.Range("L113").Value = 1500
$ws.Range("M113").Value = 988
$ws.Range("N113").Value = -5840

# GSM row 122: Awarding Academic Excellence | Ametrine
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 2309.2683
$ws.Range("I122").Value = 1953.2273
$ws.Range("J122").Value = 2721.5264
$ws.Range("K122").Value = 5859.6819
$ws.Range("L122").Value = 8164.5792
$ws.Range("M122").Value = -3409.6819
$ws.Range("N122").Value = -13064.5792

# GSM row 126: Gold Rush Order | Phrygian Gold Ingot
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 4017.9214
$ws.Range("I126").Value = 2862.9348
$ws.Range("J126").Value = 5253.4883
$ws.Range("K126").Value = 8588.804400000001
$ws.Range("L126").Value = 15760.4649
$ws.Range("M126").Value = -6118.804400000001
$ws.Range("N126").Value = -20700.4649

# GSM row 132: On Board for Lar | Lar Ingot
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 4312
$ws.Range("I132").Value = 1518.3334
$ws.Range("J132").Value = 5150.1
$ws.Range("K132").Value = 4555.0002
$ws.Range("L132").Value = 15450.3
$ws.Range("M132").Value = -2025.0002
$ws.Range("N132").Value = -20510.3

# GSM row 137: Sew Excited | Cobalt Tungsten Needle
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H137").Value = 72655.11
$ws.Range("J137").Value = 72655.11
$ws.Range("L137").Value = 72655.11
$ws.Range("N137").Value = -82855.11

# LTW row 7: Tan Before the Ban | Leather
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3341
$ws.Range("I7").Value = 1556.6154
$ws.Range("J7").Value = 5274.0835
$ws.Range("K7").Value = 1556.6154
$ws.Range("L7").Value = 5274.0835
$ws.Range("M7").Value = -1444.6154
$ws.Range("N7").Value = -5498.0835

# LTW row 40: Best Served Toad | Toad Leather
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 6604.364
$ws.Range("I40").Value = 6927.1665
$ws.Range("J40").Value = 6217
$ws.Range("K40").Value = 6927.1665
$ws.Range("L40").Value = 6217
$ws.Range("M40").Value = -6791.1665
$ws.Range("N40").Value = -6489

# LTW row 93: Hide to Go Seek | Gagana Leather
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 3714.2856
$ws.Range("I93").Value = 2000
$ws.Range("K93").Value = 2000
$ws.Range("M93").Value = -752

# LTW row 126: Battered Books | Saiga Leather
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 3341
$ws.Range("I126").Value = 1556.6154
$ws.Range("J126").Value = 5274.0835
$ws.Range("K126").Value = 4669.8462
$ws.Range("L126").Value = 15822.2505
$ws.Range("M126").Value = -2199.8462
$ws.Range("N126").Value = -20762.2505

# LTW row 127: Loyal Turncoat | Saigaskin Coat of Fending
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H127").Value = 26777.37
$ws.Range("J127").Value = 26777.37
$ws.Range("L127").Value = 26777.37
$ws.Range("N127").Value = -36697.37

# LTW row 136: Respect for Br'aax | Br'aax Leather
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 3950.5588
$ws.Range("I136").Value = 1877.0952
$ws.Range("J136").Value = 7300
$ws.Range("K136").Value = 5631.2856
$ws.Range("L136").Value = 21900
$ws.Range("M136").Value = -3081.2856
$ws.Range("N136").Value = -27000

# WVR row 62: Pride Up in Smoke | Rainbow Cloth
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 29479388
$ws.Range("I62").Value = 71431430
$ws.Range("J62").Value = 112960.3
$ws.Range("K62").Value = 71431430
$ws.Range("L62").Value = 112960.3
$ws.Range("M62").Value = -71430806
$ws.Range("N62").Value = -114208.3

# WVR row 65: Desperate for Diversionaries (L) | Rainbow Cloth
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H65").Value = 29479388
$ws.Range("I65").Value = 71431430
$ws.Range("J65").Value = 112960.3
$ws.Range("K65").Value = 357157150
$ws.Range("L65").Value = 564801.5
$ws.Range("M65").Value = -357154030
$ws.Range("N65").Value = -571041.5

# WVR row 107: Flax Wax | Bright Linen Yarn
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 684.3333
$ws.Range("I107").Value = 584.5
$ws.Range("J107").Value = 983.8333
$ws.Range("K107").Value = 1753.5
$ws.Range("L107").Value = 2951.4999
$ws.Range("M107").Value = 166.5
$ws.Range("N107").Value = -6791.4999

# WVR row 113: A Tender Table | Pixie Floss
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 173.75
$ws.Range("I113").Value = 186.21428
$ws.Range("J113").Value = 86.5
$ws.Range("K113").Value = 558.64284
$ws.Range("L113").Value = 259.5
$ws.Range("M113").Value = 1611.35716
$ws.Range("N113").Value = -4599.5

# WVR row 126: A Polished Purchase | Snow Linen
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 346246.97
$ws.Range("I126").Value = 2146.2222
$ws.Range("K126").Value = 6438.6666
$ws.Range("M126").Value = -3968.6666

# WVR row 132: Comfy Cabins | Snow Cotton Cloth
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 6668250
$ws.Range("I132").Value = 614.73334
$ws.Range("J132").Value = 16669702
$ws.Range("K132").Value = 1844.20002
$ws.Range("L132").Value = 50009106
$ws.Range("M132").Value = 685.79998
$ws.Range("N132").Value = -50014166

# WVR row 136: Weaving the Envelope | Sarcenet Cloth
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 3019.1
$ws.Range("I136").Value = 1689.5555
$ws.Range("J136").Value = 4106.909
$ws.Range("K136").Value = 5068.666499999999
$ws.Range("L136").Value = 12320.727
$ws.Range("M136").Value = -2518.666499999999
$ws.Range("N136").Value = -17420.727
